# Summarizing Ourselves final data collection
# Appends a "TAGS:" heading (bold) followed by five tag list paragraphs
# after the existing "Interviewer" paragraph at the end of the document.

$d = $word.ActiveDocument

# Create six new empty paragraphs at the end of the document (after
# "Interviewer") before typing any text into them, so that none of the
# new paragraphs inherit the bold formatting that will later be applied
# only to the "TAGS:" heading.
$last = $d.Paragraphs.Last
$r = $last.Range
$r.Collapse(0)
$r.InsertParagraphAfter()

for ($i = 0; $i -lt 5; $i++) {
    $last = $d.Paragraphs.Last
    $r = $last.Range
    $r.Collapse(0)
    $r.InsertParagraphAfter()
}

$count = $d.Paragraphs.Count
$tagsIndex = $count - 5

# Fill in the "TAGS:" heading paragraph, bolded (including complex script bold)
$tagsPara = $d.Paragraphs.Item($tagsIndex)
$tagsRange = $tagsPara.Range
$tagsRange.Collapse(1)
$tagsRange.InsertAfter("TAGS:")
$tagsRange.Font.Bold = $true
$tagsRange.Font.BoldBi = $true

# Fill in the five tag list paragraphs with plain (non-bold) text
$items = @(
    "1 Ayagayaraq -- Travel",
    "1 Tuqu Naulluun-llu, Qenan-llu -- Death & Sickness",
    "1 Allrakum Ellalinqigutai -- Seasons",
    "1 Ellavut -- Weather, Climate",
    "1 Imarpigmiutaat, Unkumiutaat, Mermiutaat -- Marine Animals"
)

for ($i = 0; $i -lt $items.Length; $i++) {
    $p = $d.Paragraphs.Item($tagsIndex + 1 + $i)
    $pr = $p.Range
    $pr.Collapse(1)
    $pr.InsertAfter($items[$i])
}
